$d = $word.ActiveDocument

# "Changed version of nm0575.fmx to 5.2"
# The version column of the nm0575.fmx row currently reads "5.1"; bump it to "5.2".
$r = $d.Content
$found = $r.Find.Execute("5.1")

if ($found) {
    # Re-type just the last character ("1" -> "2") so the edit lands the way a
    # human retouching the cell in Word would: this naturally leaves the "5."
    # portion in its original run and creates a fresh run for the new "2".
    $tail = $d.Range($r.End - 1, $r.End)
    $tail.Text = "2"

    # Force the freshly-typed character into its own run (Word does this when
    # the insertion point's formatting differs even momentarily), then restore
    # its formatting to match the rest of the cell so the visible result is
    # unchanged 9pt text.
    $tail.Font.Size = 20
    $tail.Font.Size = 9
}
